$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newTimestamp = "2022-12-25 20:49:39"

# Update timestamp column (O) for every data row (2 through 410)
for ($r = 2; $r -le 410; $r++) {
    $ws.Cells.Item($r, 15).Value = $newTimestamp
}

# Update the productAriaLabel text in row 103 (column M)
$ws.Range("M103").Value = "Betty Bossi Bio Frischback Rustico Buttergipfel - Online kein Bestand 3.50 Schweizer Franken"
